$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.5612598015441392
$ws.Range("C2").Value = 0.2063980941664276
$ws.Range("E2").Value = 0.1337229955632111
$ws.Range("F2").Value = 0.4443680307746263
$ws.Range("G2").Value = 0.0024990320184048
$ws.Range("K2").Value = 0.2713727032321742
$ws.Range("L2").Value = 0.1924470473464268
$ws.Range("M2").Value = 0.1614391636097032
$ws.Range("O2").Value = 4.429621943480214
$ws.Range("B3").Value = 0.5267956353129932
$ws.Range("C3").Value = 0.2067292841341768
$ws.Range("E3").Value = 0.1343082166931868
$ws.Range("F3").Value = 0.387822817061874
$ws.Range("G3").Value = 0.002501422356449482
$ws.Range("K3").Value = 0.240665292669604
$ws.Range("L3").Value = 0.1899127339680362
$ws.Range("M3").Value = 0.1549442551754225
$ws.Range("O3").Value = 4.472038509024884
$ws.Range("B4").Value = 0.5057995710549221
$ws.Range("C4").Value = 0.2069618841588827
$ws.Range("E4").Value = 0.1347221439121853
$ws.Range("F4").Value = 0.3531389305168915
$ws.Range("G4").Value = 0.00250296821329065
$ws.Range("K4").Value = 0.2218039225393653
$ws.Range("L4").Value = 0.1884471570645019
$ws.Range("M4").Value = 0.1510193880576232
$ws.Range("O4").Value = 4.500279134911096
$ws.Range("B5").Value = 0.4972856083549857
$ws.Range("C5").Value = 0.2070640573938185
$ws.Range("E5").Value = 0.1349045729691625
$ws.Range("F5").Value = 0.3390132514313251
$ws.Range("G5").Value = 0.002503617878382494
$ws.Range("K5").Value = 0.2141164810463749
$ws.Range("L5").Value = 0.1878727446995612
$ws.Range("M5").Value = 0.1494359414391369
$ws.Range("O5").Value = 4.512339807025398
$ws.Range("B6").Value = 0.4958744318662696
$ws.Range("C6").Value = 0.2070814702791992
$ws.Range("E6").Value = 0.1349356962791521
$ws.Range("F6").Value = 0.336668177824194
$ws.Range("G6").Value = 0.002503726947568756
$ws.Range("K6").Value = 0.2128399233403968
$ws.Range("L6").Value = 0.1877787442710144
$ws.Range("M6").Value = 0.149173979292307
$ws.Range("O6").Value = 4.514375835742257
$ws.Range("B7").Value = 0.5056845775084753
$ws.Range("C7").Value = 0.2069632321498496
$ws.Range("E7").Value = 0.1347245485116435
$ws.Range("F7").Value = 0.3529483938344953
$ws.Range("G7").Value = 0.002502976894706084
$ws.Range("K7").Value = 0.2217002515944699
$ws.Range("L7").Value = 0.1884393178464165
$ws.Range("M7").Value = 0.1509979683290403
$ws.Range("O7").Value = 4.500439552741142
$ws.Range("B8").Value = 0.5493426409530855
$ws.Range("C8").Value = 0.2065062371938637
$ws.Range("E8").Value = 0.1339134602371352
$ws.Range("F8").Value = 0.4248636149813336
$ws.Range("G8").Value = 0.002499840015555668
$ws.Range("K8").Value = 0.2607864816713459
$ws.Range("L8").Value = 0.1915544700209679
$ws.Range("M8").Value = 0.1591867053503933
$ws.Range("O8").Value = 4.443791263134926
$ws.Range("B9").Value = 0.6362442573989142
$ws.Range("C9").Value = 0.205840763659296
$ws.Range("E9").Value = 0.132755271812762
$ws.Range("F9").Value = 0.5661985755041457
$ws.Range("G9").Value = 0.002494306291712215
$ws.Range("K9").Value = 0.3373635914175566
$ws.Range("L9").Value = 0.1983790427813261
$ws.Range("M9").Value = 0.1757406744412506
$ws.Range("O9").Value = 4.350133393697291
$ws.Range("B10").Value = 0.7008542045780644
$ws.Range("C10").Value = 0.2054906760205171
$ws.Range("E10").Value = 0.1321668425962663
$ws.Range("F10").Value = 0.6702781546542269
$ws.Range("G10").Value = 0.002490613530813152
$ws.Range("K10").Value = 0.3935654702996771
$ws.Range("L10").Value = 0.203827056814518
$ws.Range("M10").Value = 0.1882008520717378
$ws.Range("O10").Value = 4.29194798214553
$ws.Range("B11").Value = 0.7304084605032415
$ws.Range("C11").Value = 0.2053611852232891
$ws.Range("E11").Value = 0.1319559253749674
$ws.Range("F11").Value = 0.7176906081379002
$ws.Range("G11").Value = 0.002489013776976147
$ws.Range("K11").Value = 0.4191172243212122
$ws.Range("L11").Value = 0.2063992885752555
$ws.Range("M11").Value = 0.193933184276105
$ws.Range("O11").Value = 4.26778510014023
$ws.Range("B12").Value = 0.7416228034570906
$ws.Range("C12").Value = 0.2053163997970202
$ws.Range("E12").Value = 0.1318841993093862
$ws.Range("F12").Value = 0.7356546913071611
$ws.Range("G12").Value = 0.002488419451076372
$ws.Range("K12").Value = 0.4287905057372257
$ws.Range("L12").Value = 0.2073867748942035
$ws.Range("M12").Value = 0.196112991905494
$ws.Range("O12").Value = 4.258966877865504
$ws.Range("B13").Value = 0.7392065906139464
$ws.Range("C13").Value = 0.2053258566037073
$ws.Range("E13").Value = 0.131899284903934
$ws.Range("F13").Value = 0.7317853510981394
$ws.Range("G13").Value = 0.002488546940734113
$ws.Range("K13").Value = 0.4267073155215257
$ws.Range("L13").Value = 0.2071735054273773
$ws.Range("M13").Value = 0.1956431286849565
$ws.Range("O13").Value = 4.260851283442236
$ws.Range("B14").Value = 0.7313306181643213
$ws.Range("C14").Value = 0.2053574156964828
$ws.Range("E14").Value = 0.1319498613229406
$ws.Range("F14").Value = 0.7191683204515869
$ws.Range("G14").Value = 0.002488964651915337
$ws.Range("K14").Value = 0.4199131061812977
$ws.Range("L14").Value = 0.2064802607254705
$ws.Range("M14").Value = 0.1941123367408579
$ws.Range("O14").Value = 4.267052970365256
$ws.Range("B15").Value = 0.7265093065127246
$ws.Range("C15").Value = 0.2053772991218139
$ws.Range("E15").Value = 0.1319819008386531
$ws.Range("F15").Value = 0.7114413442032514
$ws.Range("G15").Value = 0.002489222003736393
$ws.Range("K15").Value = 0.4157511019324716
$ws.Range("L15").Value = 0.2060573764154583
$ws.Range("M15").Value = 0.1931758635843721
$ws.Range("O15").Value = 4.270894888854428
$ws.Range("B16").Value = 0.6989259876410472
$ws.Range("C16").Value = 0.2054997349120953
$ws.Range("E16").Value = 0.1321817671408461
$ws.Range("F16").Value = 0.6671810134426437
$ws.Range("G16").Value = 0.002490719685703947
$ws.Range("K16").Value = 0.3918952640513567
$ws.Range("L16").Value = 0.2036608394170685
$ws.Range("M16").Value = 0.1878275105745359
$ws.Range("O16").Value = 4.293573497242534
$ws.Range("B17").Value = 0.6820457784535847
$ws.Range("C17").Value = 0.20558244734341
$ws.Range("E17").Value = 0.1323189035468584
$ws.Range("F17").Value = 0.6400460337125793
$ws.Range("G17").Value = 0.002491658943884884
$ws.Range("K17").Value = 0.3772563495589907
$ws.Range("L17").Value = 0.2022146460240606
$ws.Range("M17").Value = 0.1845628069225143
$ws.Range("O17").Value = 4.308076800883811
$ws.Range("B18").Value = 0.6723520960746043
$ws.Range("C18").Value = 0.2056328248412456
$ws.Range("E18").Value = 0.1324031246191879
$ws.Range("F18").Value = 0.6244449056556647
$ws.Range("G18").Value = 0.002492206722903628
$ws.Range("K18").Value = 0.3688350727947238
$ws.Range("L18").Value = 0.2013916763124257
$ws.Range("M18").Value = 0.1826910798483325
$ws.Range("O18").Value = 4.316635756186656
$ws.Range("B19").Value = 0.6690726404635825
$ws.Range("C19").Value = 0.2056503642293563
$ws.Range("E19").Value = 0.1324325588041138
$ws.Range("F19").Value = 0.619163680173358
$ws.Range("G19").Value = 0.00249239348857174
$ws.Range("K19").Value = 0.3659835545180101
$ws.Range("L19").Value = 0.2011145535967387
$ws.Range("M19").Value = 0.1820583871084125
$ws.Range("O19").Value = 4.319570943630865
$ws.Range("B20").Value = 0.6838411196497134
$ws.Range("C20").Value = 0.2055733525338468
$ws.Range("E20").Value = 0.1323037522101416
$ws.Range("F20").Value = 0.642933953830422
$ws.Range("G20").Value = 0.002491558177988583
$ws.Range("K20").Value = 0.3788148313460056
$ws.Range("L20").Value = 0.2023676810885462
$ws.Range("M20").Value = 0.1849097155602237
$ws.Range("O20").Value = 4.306510435723879
$ws.Range("B21").Value = 0.7336433687692647
$ws.Range("C21").Value = 0.2053480309238438
$ws.Range("E21").Value = 0.1319347849332004
$ws.Range("F21").Value = 0.7228739723491628
$ws.Range("G21").Value = 0.002488841649157161
$ws.Range("K21").Value = 0.4219088049813138
$ws.Range("L21").Value = 0.2066835193152343
$ws.Range("M21").Value = 0.1945617215217084
$ws.Range("O21").Value = 4.265222380200299
$ws.Range("B22").Value = 0.7663245492693136
$ws.Range("C22").Value = 0.2052255282035489
$ws.Range("E22").Value = 0.1317411026149671
$ws.Range("F22").Value = 0.7751780083420101
$ws.Range("G22").Value = 0.002487133050891589
$ws.Range("K22").Value = 0.45005774948865
$ws.Range("L22").Value = 0.2095824656465197
$ws.Range("M22").Value = 0.2009228368197498
$ws.Range("O22").Value = 4.240171931558052
$ws.Range("B23").Value = 0.7488700699673529
$ws.Range("C23").Value = 0.205288654804356
$ws.Range("E23").Value = 0.1318401381521888
$ws.Range("F23").Value = 0.7472568307830727
$ws.Range("G23").Value = 0.002488038865998969
$ws.Range("K23").Value = 0.4350357040290191
$ws.Range("L23").Value = 0.208028100819476
$ws.Range("M23").Value = 0.1975229841115862
$ws.Range("O23").Value = 4.25336484789247
$ws.Range("B24").Value = 0.6830294116625737
$ws.Range("C24").Value = 0.2055774554959271
$ws.Range("E24").Value = 0.1323105853716147
$ws.Range("F24").Value = 0.6416283278902171
$ws.Range("G24").Value = 0.002491603710060384
$ws.Range("K24").Value = 0.3781102577974025
$ws.Range("L24").Value = 0.2022984675603396
$ws.Range("M24").Value = 0.1847528619755963
$ws.Range("O24").Value = 4.307217901600467
$ws.Range("B25").Value = 0.6125994348638528
$ws.Range("C25").Value = 0.2059962820018768
$ws.Range("E25").Value = 0.1330224180165711
$ws.Range("F25").Value = 0.5279251897347166
$ws.Range("G25").Value = 0.002495737566298129
$ws.Range("K25").Value = 0.3166567158696694
$ws.Range("L25").Value = 0.1964564301926472
$ws.Range("M25").Value = 0.1712097480929842
$ws.Range("O25").Value = 4.373604443565625
